$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

# Update hashcode values (column B) for the rows that changed in this
# automatic hashcode refresh commit.
$ws.Range("B34").Value = "186ccd27c01fbd09715cdc21aef20178"
$ws.Range("B162").Value = "f769d9ad8868add77ed3b779af21d829"
$ws.Range("B175").Value = "11d9a85c51d17e04903382def657c744"
$ws.Range("B180").Value = "249a45a2959fddef02e8381674ad208b"
$ws.Range("B191").Value = "f6a1624b20a7c32238733f979dcbf78e"
$ws.Range("B213").Value = "26775bab7f6472923f25ffba80d02231"
$ws.Range("B338").Value = "23e0f2ebb1ed868a183939e0a30e00f3"
$ws.Range("B423").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B488").Value = "236f85fc893a94d3911ed3f3a90febe0"
$ws.Range("B516").Value = "07500c8bb13b31c526a0c97fdde3999b"
$ws.Range("B524").Value = "f9039b034685abd2f12c7a2d6a855dc4"
$ws.Range("B535").Value = "a2e40a52b2f1138e0633df4c4146ba9b"
$ws.Range("B596").Value = "2b967f8d71c14c353a5c5f0c1939ed3c"
$ws.Range("B678").Value = "c1b10d31595f882695165018f4e34e13"
$ws.Range("B716").Value = "b3dff372473427a4fd582fc76c2a3e3f"
$ws.Range("B738").Value = "696e5e106d883e1866d4e144dd54d9fa"
$ws.Range("B741").Value = "654d4770dcce46793cfce891887dfa1a"
$ws.Range("B754").Value = "67a632afcce75e9838d6a233116ab671"
$ws.Range("B828").Value = "45a541e4ef4bda24818c78b4fd4f873d"
$ws.Range("B837").Value = "f1a4971618a3340bd5f04bd6b09c480d"
